# Updated The Mapping Sheets to Some Point
#
# 1. Add two new worksheets ("BankAccount", "Mobile Wallets") at the end of
#    the workbook describing the bank_accounts / mobile_wallets MySQL tables.
# 2. Update the remembered cell selection on a couple of the pre-existing
#    sheets ("withdrow requests", "transactions table").
# 3. Leave "Mobile Wallets" as the active/selected sheet+cell, matching the
#    author's last-saved UI state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1a. "BankAccount" sheet - appended after the existing "transactions table"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bank = $wb.Worksheets.Add($null, $lastSheet)
$bank.Name = "BankAccount"

# Data entry order matters for shared-string allocation: the author typed
# the SQL column/type/constraint tokens (columns B:G) down every row first,
# then went back and filled in the header row and the repeated "Table Name"
# column (A) - so we reproduce cells in that same order.
$bank.Cells.Item(2,2).Value = '`id`'
$bank.Cells.Item(2,3).Value = 'bigint'
$bank.Cells.Item(2,4).Value = 'unsigned'
$bank.Cells.Item(2,5).Value = 'NOT'
$bank.Cells.Item(2,6).Value = 'NULL'
$bank.Cells.Item(2,7).Value = 'AUTO_INCREMENT,'

$bank.Cells.Item(3,2).Value = '`user_id`'
$bank.Cells.Item(3,3).Value = 'bigint'
$bank.Cells.Item(3,4).Value = 'NOT'
$bank.Cells.Item(3,5).Value = 'NULL,'

$bank.Cells.Item(4,2).Value = '`type`'
$bank.Cells.Item(4,3).Value = 'enum(''WITHDRAWAL'',''CASHOUT'')'
$bank.Cells.Item(4,4).Value = 'DEFAULT'
$bank.Cells.Item(4,5).Value = 'NULL,'

$bank.Cells.Item(5,2).Value = '`bank_name`'
$bank.Cells.Item(5,3).Value = 'varchar(191)'
$bank.Cells.Item(5,4).Value = 'NOT'
$bank.Cells.Item(5,5).Value = 'NULL,'

$bank.Cells.Item(6,2).Value = '`bank_address`'
$bank.Cells.Item(6,3).Value = 'text'
$bank.Cells.Item(6,4).Value = 'NOT'
$bank.Cells.Item(6,5).Value = 'NULL,'

$bank.Cells.Item(7,2).Value = '`number`'
$bank.Cells.Item(7,3).Value = 'bigint'
$bank.Cells.Item(7,4).Value = 'NOT'
$bank.Cells.Item(7,5).Value = 'NULL,'

$bank.Cells.Item(8,2).Value = '`name`'
$bank.Cells.Item(8,3).Value = 'varchar(191)'
$bank.Cells.Item(8,4).Value = 'DEFAULT'
$bank.Cells.Item(8,5).Value = 'NULL,'

$bank.Cells.Item(9,2).Value = '`holder_name`'
$bank.Cells.Item(9,3).Value = 'varchar(191)'
$bank.Cells.Item(9,4).Value = 'NOT'
$bank.Cells.Item(9,5).Value = 'NULL,'

$bank.Cells.Item(10,2).Value = '`swift_code`'
$bank.Cells.Item(10,3).Value = 'varchar(191)'
$bank.Cells.Item(10,4).Value = 'NOT'
$bank.Cells.Item(10,5).Value = 'NULL,'

$bank.Cells.Item(11,2).Value = '`created_at`'
$bank.Cells.Item(11,3).Value = 'timestamp'
$bank.Cells.Item(11,4).Value = 'NULL'
$bank.Cells.Item(11,5).Value = 'DEFAULT'
$bank.Cells.Item(11,6).Value = 'NULL,'

$bank.Cells.Item(12,2).Value = '`updated_at`'
$bank.Cells.Item(12,3).Value = 'timestamp'
$bank.Cells.Item(12,4).Value = 'NULL'
$bank.Cells.Item(12,5).Value = 'DEFAULT'
$bank.Cells.Item(12,6).Value = 'NULL,'

$bank.Cells.Item(13,2).Value = '`deleted_at`'
$bank.Cells.Item(13,3).Value = 'timestamp'
$bank.Cells.Item(13,4).Value = 'NULL'
$bank.Cells.Item(13,5).Value = 'DEFAULT'
$bank.Cells.Item(13,6).Value = 'NULL,'

# header row
$bank.Cells.Item(1,1).Value = 'Table Name'
$bank.Cells.Item(1,2).Value = 'Column'
$bank.Cells.Item(1,3).Value = 'DataType'

# table-name column, same value repeated down every data row
for ($r = 2; $r -le 13; $r++) {
    $bank.Cells.Item($r,1).Value = 'bank_accounts'
}

$bank.Columns.Item(1).ColumnWidth = 48.88671875
$bank.Columns.Item(3).ColumnWidth = 29.109375

# ---------------------------------------------------------------------
# 1b. "Mobile Wallets" sheet - appended after "BankAccount"
# ---------------------------------------------------------------------
$wallet = $wb.Worksheets.Add($null, $bank)
$wallet.Name = "Mobile Wallets"

$wallet.Cells.Item(2,2).Value = '`id`'
$wallet.Cells.Item(2,3).Value = 'bigint'
$wallet.Cells.Item(2,4).Value = 'unsigned'
$wallet.Cells.Item(2,5).Value = 'NOT'
$wallet.Cells.Item(2,6).Value = 'NULL'
$wallet.Cells.Item(2,7).Value = 'AUTO_INCREMENT,'

$wallet.Cells.Item(3,2).Value = '`user_id`'
$wallet.Cells.Item(3,3).Value = 'int'
$wallet.Cells.Item(3,4).Value = 'unsigned'
$wallet.Cells.Item(3,5).Value = 'NOT'
$wallet.Cells.Item(3,6).Value = 'NULL,'

$wallet.Cells.Item(4,2).Value = '`holder_name`'
$wallet.Cells.Item(4,3).Value = 'varchar(191)'
$wallet.Cells.Item(4,4).Value = 'DEFAULT'
$wallet.Cells.Item(4,5).Value = 'NULL,'

$wallet.Cells.Item(5,2).Value = '`mobile_number`'
$wallet.Cells.Item(5,3).Value = 'varchar(191)'
$wallet.Cells.Item(5,4).Value = 'NOT'
$wallet.Cells.Item(5,5).Value = 'NULL,'

$wallet.Cells.Item(6,2).Value = '`created_at`'
$wallet.Cells.Item(6,3).Value = 'timestamp'
$wallet.Cells.Item(6,4).Value = 'NULL'
$wallet.Cells.Item(6,5).Value = 'DEFAULT'
$wallet.Cells.Item(6,6).Value = 'NULL,'

$wallet.Cells.Item(7,2).Value = '`updated_at`'
$wallet.Cells.Item(7,3).Value = 'timestamp'
$wallet.Cells.Item(7,4).Value = 'NULL'
$wallet.Cells.Item(7,5).Value = 'DEFAULT'
$wallet.Cells.Item(7,6).Value = 'NULL,'

$wallet.Cells.Item(8,2).Value = '`deleted_at`'
$wallet.Cells.Item(8,3).Value = 'timestamp'
$wallet.Cells.Item(8,4).Value = 'NULL'
$wallet.Cells.Item(8,5).Value = 'DEFAULT'
$wallet.Cells.Item(8,6).Value = 'NULL,'

# header row
$wallet.Cells.Item(1,1).Value = 'Table Name'
$wallet.Cells.Item(1,2).Value = 'Column Name'
$wallet.Cells.Item(1,3).Value = 'DataType'

# table-name column, same value repeated down every data row
for ($r = 2; $r -le 8; $r++) {
    $wallet.Cells.Item($r,1).Value = 'mobile_wallets'
}

$wallet.Columns.Item(1).ColumnWidth = 13.44140625
$wallet.Columns.Item(2).ColumnWidth = 15
$wallet.Columns.Item(3).ColumnWidth = 11.21875

# ---------------------------------------------------------------------
# 2. Refresh remembered selections on the pre-existing sheets that moved
# ---------------------------------------------------------------------
$withdrow = $wb.Worksheets.Item("withdrow requests")
$withdrow.Activate()
$withdrow.Range("A4").Select()

$transactions = $wb.Worksheets.Item("transactions table")
$transactions.Activate()
$transactions.Range("B5").Select()

# ---------------------------------------------------------------------
# 3. Leave the new sheets selected as the author last left them, with
#    "Mobile Wallets" as the final active tab.
# ---------------------------------------------------------------------
$bank.Activate()
$bank.Range("E23").Select()

$wallet.Activate()
$wallet.Range("J10").Select()
